$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 42; this shifts existing rows 42-60 down to 43-61
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record
$ws.Cells.Item(42, 1).Value = 2
$ws.Cells.Item(42, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = 44755
$ws.Cells.Item(42, 4).NumberFormat = $ws.Cells.Item(43, 4).NumberFormat
$ws.Cells.Item(42, 5).Value = 4
$ws.Cells.Item(42, 6).Value = 100112022
$ws.Cells.Item(42, 7).Value = "Arveja Verde"
$ws.Cells.Item(42, 8).Value = "Perfection"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 200
$ws.Cells.Item(42, 11).Value = 30000
$ws.Cells.Item(42, 12).Value = 32000
$ws.Cells.Item(42, 13).Value = 31000
$ws.Cells.Item(42, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 1240
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
